# Generate Report for Handback
# Updates the localization status report to reflect that the
# 9ba2b223-8801-4341-9a8c-da5d7a438d89.md file has been handed back
# (in sync with en-US) for both the zh-cn and de-de locales.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = $statusHandedBack
$wsZhCn.Range("G3").Value = "2016-03-09 14:16:13"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = $statusHandedBack
$wsDeDe.Range("G3").Value = "2016-03-09 14:16:18"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack
